# "fixed bug in calculation of q-coefficients"
#
# The sheet was missing a "time spent [h]" column (the q-coefficient that
# turns minutes into hours) and one working day (2014-02-21, 18:30-22:00)
# had never been added to the table, so it was missing from every sum.
# This adds column G with that conversion for every data row, adds the
# missing row, and shifts the three summary rows ("sum [min]", "sum [h]",
# "sum [working weeks]") down by one so they keep covering the full range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 14 used to be an empty trailer row; it is now real data ------
$ws.Range("A14").Value = 2014
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 0.77083333333333337
$ws.Range("E14").Value = 0.91666666666666663
$ws.Range("D14").NumberFormat = "hh:mm;@"
$ws.Range("E14").NumberFormat = "hh:mm;@"
$ws.Range("F14").Formula = "=(E14-D14)*24*60"
$ws.Range("F14").NumberFormat = "0"

# --- new column header and q-coefficient ("time spent [h]") column ----
$ws.Range("G1").Value = "time spent [h]"
$ws.Range("G2").Formula = "=F2/60"
$ws.Range("G3:G14").Formula = "=F3/60"
$ws.Range("G1:G14").NumberFormat = "0.00"

# --- the blank trailer row moves from row 14 to row 15 -----------------
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("E15").Style = "Normal"
$ws.Range("D15").NumberFormat = "hh:mm;@"
$ws.Range("E15").NumberFormat = "hh:mm;@"
$ws.Range("F15").NumberFormat = "0"

# --- "sum [min]" moves from row 15 to row 16, now covering F2:F15 -----
$ws.Range("E16").Value = "sum [min]"
$ws.Range("E16").HorizontalAlignment = -4152
$ws.Range("F16").Formula = "=SUM(F2:F15)"
$ws.Range("F16").NumberFormat = "0"

# --- "sum [h]" moves from row 16 to row 17 -----------------------------
$ws.Range("E17").Value = "sum [h]"
$ws.Range("E17").HorizontalAlignment = -4152
$ws.Range("F17").Formula = "=F16/60"
$ws.Range("F17").NumberFormat = "0.00"

# --- "sum [working weeks]" moves from row 17 to row 18 -----------------
$ws.Range("E18").Value = "sum [working weeks]"
$ws.Range("E18").HorizontalAlignment = -4152
$ws.Range("F18").Formula = "=F17/38.5"
$ws.Range("F18").NumberFormat = "0.00"

# --- column width for the new column G (closest reachable to 13.71) ---
$ws.Columns("G").ColumnWidth = 12.8

# --- selection matches the post-edit cursor position --------------------
$ws.Range("I14").Select()
